# Update gh-pages data output (想去人数 counts refreshed) for the
# 上海-漫展信息 workbook. Column F on every sheet holds the "想去人数"
# (number of people interested) counter that was re-scraped; only those
# values changed between the previous and the newly generated snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 312
$ws.Range("F3").Value  = 1089
$ws.Range("F4").Value  = 1236
$ws.Range("F5").Value  = 1112
$ws.Range("F6").Value  = 3332
$ws.Range("F7").Value  = 57
$ws.Range("F9").Value  = 1165
$ws.Range("F10").Value = 738
$ws.Range("F11").Value = 577
$ws.Range("F12").Value = 279
$ws.Range("F13").Value = 51
$ws.Range("F14").Value = 136
$ws.Range("F15").Value = 646
$ws.Range("F16").Value = 1698
$ws.Range("F17").Value = 1698
$ws.Range("F18").Value = 35
$ws.Range("F19").Value = 334
$ws.Range("F20").Value = 19
$ws.Range("F21").Value = 42
$ws.Range("F22").Value = 631
$ws.Range("F23").Value = 378
$ws.Range("F25").Value = 664
$ws.Range("F26").Value = 77719
$ws.Range("F27").Value = 77719
$ws.Range("F29").Value = 651
$ws.Range("F30").Value = 33377
$ws.Range("F31").Value = 33377
$ws.Range("F32").Value = 472
$ws.Range("F33").Value = 16
$ws.Range("F34").Value = 13
$ws.Range("F35").Value = 45
$ws.Range("F36").Value = 17
$ws.Range("F37").Value = 941
$ws.Range("F38").Value = 264
$ws.Range("F39").Value = 153
$ws.Range("F40").Value = 537
$ws.Range("F42").Value = 1167
$ws.Range("F43").Value = 5415
$ws.Range("F44").Value = 745
$ws.Range("F45").Value = 437
$ws.Range("F47").Value = 357
$ws.Range("F50").Value = 17
$ws.Range("F51").Value = 41
$ws.Range("F52").Value = 5

# ---- Sheet: 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 21
$ws.Range("F15").Value = 1110
$ws.Range("F17").Value = 71
$ws.Range("F18").Value = 405
$ws.Range("F20").Value = 64
$ws.Range("F23").Value = 494
$ws.Range("F24").Value = 10
$ws.Range("F35").Value = 1648
$ws.Range("F42").Value = 18
$ws.Range("F43").Value = 29
$ws.Range("F46").Value = 812
$ws.Range("F47").Value = 102
$ws.Range("F48").Value = 102
$ws.Range("F49").Value = 40

# ---- Sheet: 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 707
$ws.Range("F5").Value = 548
$ws.Range("F6").Value = 566

# ---- Sheet: 全部类型 (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 707
$ws.Range("F4").Value  = 312
$ws.Range("F5").Value  = 548
$ws.Range("F6").Value  = 1089
$ws.Range("F7").Value  = 1236
$ws.Range("F9").Value  = 1112
$ws.Range("F10").Value = 3332
$ws.Range("F11").Value = 57
$ws.Range("F14").Value = 1165
$ws.Range("F15").Value = 738
$ws.Range("F18").Value = 566
$ws.Range("F20").Value = 577
$ws.Range("F22").Value = 51
$ws.Range("F23").Value = 1698
$ws.Range("F24").Value = 1698
$ws.Range("F25").Value = 35
$ws.Range("F26").Value = 334
$ws.Range("F28").Value = 19
$ws.Range("F29").Value = 42
$ws.Range("F30").Value = 631
$ws.Range("F31").Value = 378
$ws.Range("F32").Value = 664
$ws.Range("F33").Value = 77721
$ws.Range("F34").Value = 651
$ws.Range("F35").Value = 33377
$ws.Range("F36").Value = 472
$ws.Range("F37").Value = 16
$ws.Range("F38").Value = 13
$ws.Range("F39").Value = 45
$ws.Range("F40").Value = 17
$ws.Range("F42").Value = 264
$ws.Range("F44").Value = 537
$ws.Range("F46").Value = 5415
$ws.Range("F47").Value = 1648
$ws.Range("F53").Value = 29
$ws.Range("F55").Value = 102
